$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "iaguilera"
$ws.Range("C4").Value = "Ignacio Aguilera"
$ws.Range("B4").Value = '$2b$10$NKSgcAXRuxSqKkECocB2/egUEGcVN02XkKrVpbmFo3js.gGmrrtQu'
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D1").Value = "alias"
$ws.Range("D3").Value = "Maicita"
$ws.Range("D2").Value = "Administrador"
$ws.Range("D4").Value = "Administrador"

$ws.Range("D4").Select()
